# TC for order added
# Re-works the "Order Confirmation" block (rows 51-57) on the RTM_ALL sheet
# so it holds seven distinct FR_ORDER_xx / TS_ORDER_xx rows (instead of the
# old six rows duplicated twice), then slides the Non-Functional block up
# underneath it (rows 58-64) and drops the now-empty trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM_ALL")

$ldquo = [char]0x201C
$rdquo = [char]0x201D

# --- Order Confirmation rows (51-57) ------------------------------------
# (values are written in the same order the original authoring session
# introduced them, so new shared-string entries land in the same order)
$ws.Range("A51").Value = "FR_ORDER_01"
$ws.Range("C51").Value = "Order Confirmation"
$ws.Range("D51").Value = "TS_ORDER_01"

$ws.Range("A52").Value = "FR_ORDER_02"
$ws.Range("C52").Value = "Order Confirmation"
$ws.Range("D52").Value = "TS_ORDER_02"

$ws.Range("A53").Value = "FR_ORDER_03"
$ws.Range("C53").Value = "Order Confirmation"
$ws.Range("D53").Value = "TS_ORDER_03"

$ws.Range("A54").Value = "FR_ORDER_04"
$ws.Range("C54").Value = "Order Confirmation"
$ws.Range("D54").Value = "TS_ORDER_04"

$ws.Range("A55").Value = "FR_ORDER_05"
$ws.Range("C55").Value = "Order Confirmation"
$ws.Range("D55").Value = "TS_ORDER_05"

$ws.Range("A56").Value = "FR_ORDER_06"
$ws.Range("C56").Value = "Order Confirmation"
$ws.Range("D56").Value = "TS_ORDER_06"

$ws.Range("A57").Value = "FR_ORDER_07"
$ws.Range("C57").Value = "Order Confirmation"
$ws.Range("D57").ClearContents()

$ws.Range("B53").Value = "Order details in display"
$ws.Range("B54").Value = "Payment status should be shown"
$ws.Range("B55").Value = "Order confirmation message should be displayed"
$ws.Range("B56").Value = "Order details should be available in " + $ldquo + "My Orders" + $rdquo
$ws.Range("B57").Value = "Confirmation email/SMS should be sent"

$ws.Range("B51").Value = "Order confirmation message"
$ws.Range("B52").Value = "Unique order ID generation"

# --- Non-Functional rows slide up from 63-69 to 58-64 -------------------
$ws.Range("A58").Value = "NFR_01"
$ws.Range("B58").Value = "Pages load within 3 seconds"
$ws.Range("C58").Value = "Non-Functional"
$ws.Range("D58").Value = "TS_NFR_01"

$ws.Range("A59").Value = "NFR_02"
$ws.Range("B59").Value = "User-friendly UI"
$ws.Range("C59").Value = "Non-Functional"
$ws.Range("D59").Value = "TS_NFR_02"

$ws.Range("A60").Value = "NFR_03"
$ws.Range("B60").Value = "Meaningful error messages"
$ws.Range("C60").Value = "Non-Functional"
$ws.Range("D60").Value = "TS_NFR_03"

$ws.Range("A61").Value = "NFR_04"
$ws.Range("B61").Value = "Password masking & secure data"
$ws.Range("C61").Value = "Non-Functional"
$ws.Range("D61").Value = "TS_NFR_04"

$ws.Range("A62").Value = "NFR_05"
$ws.Range("B62").Value = "Restrict unauthorized access"
$ws.Range("C62").Value = "Non-Functional"
$ws.Range("D62").Value = "TS_NFR_05"

$ws.Range("A63").Value = "NFR_06"
$ws.Range("B63").Value = "Browser compatibility"
$ws.Range("C63").Value = "Non-Functional"
$ws.Range("D63").Value = "TS_NFR_06"

$ws.Range("A64").Value = "NFR_07"
$ws.Range("B64").Value = "Responsive design"
$ws.Range("C64").Value = "Non-Functional"
$ws.Range("D64").Value = "TS_NFR_07"

# --- Drop the old trailing rows (previously NFR_03..NFR_07 duplicates) --
$ws.Range("A65:D69").ClearContents()

# --- Restore the view's current selection --------------------------------
$ws.Range("E52:E54").Select()

Write-Host "edit complete"
